$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$ws.Rows.Item(21).Delete()
$ws.Rows.Item(1).Delete()
$ws.Rows.Item(24).Select() | Out-Null
